# The select ("Status") value wasn't being kept on edit - fill in "Done"
# for the rows where it had been dropped (C20, C22, C24, C25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = "Done"
$ws.Range("C22").Value = "Done"
$ws.Range("C24").Value = "Done"
$ws.Range("C25").Value = "Done"

# Mirror the author's in-progress cursor/scroll position: scrolled down so
# row 13 is the top visible row, with C24 as the active/selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select()
